$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 becomes a "continuation" row of the block started at row 9:
# its border style switches from s=4/5 to s=6/7 (same family as rows 3/6/8).
$ws.Range("A3:E3").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

# Row 11: a new single-line block, styled like row 4 (s=8 / s=9).
$ws.Range("A4:E4").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Rows(11).RowHeight = 43.2

$ws.Range("C11").Value = " Hey! Long time no see!"
$ws.Range("A11").Value = "SCRIPT/P01P04A/um2103.ssb"
$ws.Range("D11").Value = " Привет! Давненько не виделись!"
$ws.Range("E11").Value = " Ðñéâåó! Äàâîåîûëï îå âéäåìéòû!"
$ws.Range("B11").Value = 248

# Rows 12-13: a new two-line block, styled like row 9 (s=4 / s=5).
# Row 13 has no cell in column A (mirrors rows 3/6/8, where only the first
# row of a block carries the filename).
$ws.Range("A9:E9").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Rows(12).RowHeight = 43.2
$ws.Range("B9:E9").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122)

$ws.Range("C12").Value = " [hero], I\'m jealous of you.\nYou\'ve got a partner."
$ws.Range("C13").Value = " I want a partner, too…"
$ws.Range("A12").Value = "SCRIPT/T01P01A/m22a0501.ssb"
$ws.Range("D12").Value = " [hero], я тебе завидую.\nУ тебя есть спутник."
$ws.Range("D13").Value = " Я тоже хочу спутника..."
$ws.Range("E12").Value = " [hero], ÿ óåáå èàâéäôý.\nÔ óåáÿ åòóû òðôóîéë."
$ws.Range("E13").Value = " Ÿ óïçå öïœô òðôóîéëà..."
$ws.Range("B12").Value = 226
$ws.Range("B13").Value = 229
